# Natmi following Dr Hou advice
#
# Update the LR-pairs sheet (Tnfsf18-Tnfrsf18) to add the "FAPs" cluster:
#   - Rows 2-4 (existing Sending cluster "M2") get their numeric columns
#     refreshed (ligand/receptor-expressing cell counts, rates, weights,
#     specificities, etc.) for the new 3-cluster comparison set, and a new
#     Target cluster "FAPs" row is appended.
#   - Rows 5 (new) extend the "M2" sending-cluster block with a "sCs"
#     target-cluster row.
#   - Rows 6-9 (new) add a whole new Sending cluster "FAPs" block mirroring
#     the same four Target clusters (ECs, FAPs, M2, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the 16 numeric metric columns.
$rows = @(
    ,@("FAPs", "Tnfsf18", "Tnfrsf18", "ECs",  2, 0.6666666666666666, 0.7506253333333333, 2.251876,           0.5928186022610482, 0.5928186022610483, 3, 1,                  1.041016333333333,  3.123049, 0.3352332808110565, 0.3352332808110565, 0.7814132322137777, 7.032719089923999, 0.198732524961796,   0.198732524961796)
    ,@("FAPs", "Tnfsf18", "Tnfrsf18", "FAPs", 2, 0.6666666666666666, 0.7506253333333333, 2.251876,           0.5928186022610482, 0.5928186022610483, 2, 0.6666666666666666, 0.8548503333333333, 2.564551, 0.2752831753639714, 0.2752831753639715, 0.6416723164084444, 5.775050847675999, 0.1631929872452526,  0.1631929872452526)
    ,@("FAPs", "Tnfsf18", "Tnfrsf18", "M2",   2, 0.6666666666666666, 0.7506253333333333, 2.251876,           0.5928186022610482, 0.5928186022610483, 3, 1,                  0.4691323333333333, 1.407397, 0.1510723378703435, 0.1510723378703435, 0.3521426140857777, 3.169283526772,    0.08955849217660584, 0.08955849217660589)
    ,@("FAPs", "Tnfsf18", "Tnfrsf18", "sCs",  2, 0.6666666666666666, 0.7506253333333333, 2.251876,           0.5928186022610482, 0.5928186022610483, 3, 1,                  0.74035,             2.22105,  0.2384112059546286, 0.2384112059546286, 0.5557254655333332, 5.001529189799999, 0.1413345978773938,  0.1413345978773939)
    ,@("M2",   "Tnfsf18", "Tnfrsf18", "ECs",  3, 1,                  0.515572,            1.546716,           0.4071813977389517, 0.4071813977389517, 3, 1,                  1.041016333333333,  3.123049, 0.3352332808110565, 0.3352332808110565, 0.5367188730093334, 4.830469857084,    0.1365007558492605,  0.1365007558492605)
    ,@("M2",   "Tnfsf18", "Tnfrsf18", "FAPs", 3, 1,                  0.515572,            1.546716,           0.4071813977389517, 0.4071813977389517, 2, 0.6666666666666666, 0.8548503333333333, 2.564551, 0.2752831753639714, 0.2752831753639715, 0.4407368960573333, 3.966632064516,    0.1120901881187188,  0.1120901881187189)
    ,@("M2",   "Tnfsf18", "Tnfrsf18", "M2",   3, 1,                  0.515572,            1.546716,           0.4071813977389517, 0.4071813977389517, 3, 1,                  0.4691323333333333, 1.407397, 0.1510723378703435, 0.1510723378703435, 0.2418714953613333, 2.176843458252,    0.06151384569373763, 0.06151384569373765)
    ,@("M2",   "Tnfsf18", "Tnfrsf18", "sCs",  3, 1,                  0.515572,            1.546716,           0.4071813977389517, 0.4071813977389517, 3, 1,                  0.74035,             2.22105,  0.2384112059546286, 0.2384112059546286, 0.3817037302,       3.4353335718,      0.09707660807723476, 0.09707660807723477)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
